# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-12-24 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-25 Monday", 2)

# Update the division problems in the table. The table has 20 rows x 5
# columns; only every 4th row (1, 5, 9, 13, 17 in 1-based indexing)
# actually holds data, the others are blank spacer rows.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "22÷5="
$t.Cell(1, 2).Range.Text  = "71÷8="
$t.Cell(1, 3).Range.Text  = "17÷4="
$t.Cell(1, 4).Range.Text  = "83÷8="
$t.Cell(1, 5).Range.Text  = "49÷5="

$t.Cell(5, 1).Range.Text  = "44÷4="
$t.Cell(5, 2).Range.Text  = "75÷5="
$t.Cell(5, 3).Range.Text  = "24÷6="
$t.Cell(5, 4).Range.Text  = "93÷4="
$t.Cell(5, 5).Range.Text  = "37÷3="

$t.Cell(9, 1).Range.Text  = "72÷4="
$t.Cell(9, 2).Range.Text  = "48÷9="
$t.Cell(9, 3).Range.Text  = "40÷9="
$t.Cell(9, 4).Range.Text  = "40÷4="
$t.Cell(9, 5).Range.Text  = "19÷7="

$t.Cell(13, 1).Range.Text = "56÷4="
$t.Cell(13, 2).Range.Text = "96÷6="
$t.Cell(13, 3).Range.Text = "80÷5="
$t.Cell(13, 4).Range.Text = "35÷6="
$t.Cell(13, 5).Range.Text = "49÷6="

$t.Cell(17, 1).Range.Text = "42÷9="
$t.Cell(17, 2).Range.Text = "44÷7="
$t.Cell(17, 3).Range.Text = "84÷9="
$t.Cell(17, 4).Range.Text = "29÷2="
$t.Cell(17, 5).Range.Text = "77÷7="
